# Trade #38 closed at 2026-02-16 22:55:54 - base_strategy DOWN +0.000%
#
# Appends a new trade row (row 39) to both the "All Trades" and the
# "base_strategy" worksheets, mirroring the layout of the existing rows.

$wb = $excel.ActiveWorkbook

$targetSheets = @("All Trades", "base_strategy")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 39

    $ws.Range("A$row").Value = 38
    # Leading apostrophe forces the date string to stay plain text instead
    # of being auto-parsed as a date serial number, matching the other
    # rows above (the time-of-day string does not need this). Re-applying
    # the Normal style afterwards drops the quote-prefix formatting flag
    # so the cell keeps the plain, unstyled look of its neighbours.
    $ws.Range("B$row").Value = "'2026-02-16"
    $ws.Range("B$row").Style = "Normal"
    $ws.Range("C$row").Value = "22:55:54"
    $ws.Range("D$row").Value = "base_strategy"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("F$row").Value = 49.999998
    $ws.Range("G$row").Value = ""
    $ws.Range("H$row").Value = "OPEN"
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 100
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = ""
    $ws.Range("Q$row").Value = 0
}
